$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.309.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.097.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.095.78"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("E9").Value = "  +4.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "

$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("E12").Value = "  +4.17%  "

$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.626.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000163"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.366.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.111.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "359.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.48%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0861"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.97%  "

$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("E38").Value = "  +2.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0666"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.483.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.693"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.134.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.976"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.735"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.29%  "

